# Update computed results (pl_mw) for the 380 kV case: columns B-D and F-N
# for rows 2-25 receive their recalculated values. Columns A, E and O are
# left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.8033754279506979
$ws.Cells.Item(2, 3).Value = 0.04961386267486034
$ws.Cells.Item(2, 4).Value = 0.09646837582005929
$ws.Cells.Item(2, 6).Value = 2.287199385550934
$ws.Cells.Item(2, 7).Value = 1.601005413721211
$ws.Cells.Item(2, 8).Value = 1.451213742430042
$ws.Cells.Item(2, 9).Value = 1.454311739783506
$ws.Cells.Item(2, 10).Value = 0.1947768401313876
$ws.Cells.Item(2, 11).Value = 0.5270933037425038
$ws.Cells.Item(2, 12).Value = 0.3568937670666017
$ws.Cells.Item(2, 13).Value = 0.2595890386159851
$ws.Cells.Item(2, 14).Value = 2.698998929581851

$ws.Cells.Item(3, 2).Value = 0.771592486386993
$ws.Cells.Item(3, 3).Value = 0.04609251484683341
$ws.Cells.Item(3, 4).Value = 0.09559808043853479
$ws.Cells.Item(3, 6).Value = 2.293205292164728
$ws.Cells.Item(3, 7).Value = 1.605105955281132
$ws.Cells.Item(3, 8).Value = 1.457746528855481
$ws.Cells.Item(3, 9).Value = 1.461519784216243
$ws.Cells.Item(3, 10).Value = 0.1957119473707962
$ws.Cells.Item(3, 11).Value = 0.4943590344040842
$ws.Cells.Item(3, 12).Value = 0.3551240325136718
$ws.Cells.Item(3, 13).Value = 0.2538969803960427
$ws.Cells.Item(3, 14).Value = 2.721095116021061

$ws.Cells.Item(4, 2).Value = 0.7524263728961103
$ws.Cells.Item(4, 3).Value = 0.04390837404945103
$ws.Cells.Item(4, 4).Value = 0.0950923774158241
$ws.Cells.Item(4, 6).Value = 2.297816831157348
$ws.Cells.Item(4, 7).Value = 1.608353550521556
$ws.Cells.Item(4, 8).Value = 1.462259180989435
$ws.Cells.Item(4, 9).Value = 1.466491084947847
$ws.Cells.Item(4, 10).Value = 0.1963239729687842
$ws.Cells.Item(4, 11).Value = 0.4744775348043788
$ws.Cells.Item(4, 12).Value = 0.3541679334297925
$ws.Cells.Item(4, 13).Value = 0.2505130471135786
$ws.Cells.Item(4, 14).Value = 2.735375473536056

$ws.Cells.Item(5, 2).Value = 0.7447042991474007
$ws.Cells.Item(5, 3).Value = 0.04301277338691989
$ws.Cells.Item(5, 4).Value = 0.09489354973996456
$ws.Cells.Item(5, 6).Value = 2.299928662273047
$ws.Cells.Item(5, 7).Value = 1.609860607456881
$ws.Cells.Item(5, 8).Value = 1.464224386510622
$ws.Cells.Item(5, 9).Value = 1.468654237275899
$ws.Cells.Item(5, 10).Value = 0.1965829187659232
$ws.Cells.Item(5, 11).Value = 0.4664307831648102
$ws.Cells.Item(5, 12).Value = 0.353811243487435
$ws.Cells.Item(5, 13).Value = 0.249162117629119
$ws.Cells.Item(5, 14).Value = 2.741374244623774

$ws.Cells.Item(6, 2).Value = 0.7434274032256383
$ws.Cells.Item(6, 3).Value = 0.04286372438755848
$ws.Cells.Item(6, 4).Value = 0.09486097373365965
$ws.Cells.Item(6, 6).Value = 2.300293386209219
$ws.Cells.Item(6, 7).Value = 1.610121947497348
$ws.Cells.Item(6, 8).Value = 1.464558338036909
$ws.Cells.Item(6, 9).Value = 1.469021723861104
$ws.Cells.Item(6, 10).Value = 0.196626493284942
$ws.Cells.Item(6, 11).Value = 0.4650979681753142
$ws.Cells.Item(6, 12).Value = 0.3537540073147269
$ws.Cells.Item(6, 13).Value = 0.2489394950085035
$ws.Cells.Item(6, 14).Value = 2.742381173894909

$ws.Cells.Item(7, 2).Value = 0.7523218721257479
$ws.Cells.Item(7, 3).Value = 0.04389631813058514
$ws.Cells.Item(7, 4).Value = 0.09508966654129836
$ws.Cells.Item(7, 6).Value = 2.29784436993608
$ws.Cells.Item(7, 7).Value = 1.608373131540972
$ws.Cells.Item(7, 8).Value = 1.462285172995124
$ws.Cells.Item(7, 9).Value = 1.466519701856917
$ws.Cells.Item(7, 10).Value = 0.1963274265487787
$ws.Cells.Item(7, 11).Value = 0.4743687897382642
$ws.Cells.Item(7, 12).Value = 0.3541629895172989
$ws.Cells.Item(7, 13).Value = 0.2504947142494345
$ws.Cells.Item(7, 14).Value = 2.735455648551909

$ws.Cells.Item(8, 2).Value = 0.7923446244064678
$ws.Cells.Item(8, 3).Value = 0.04840426585566604
$ws.Cells.Item(8, 4).Value = 0.09616237463490762
$ws.Cells.Item(8, 6).Value = 2.289078629224299
$ws.Cells.Item(8, 7).Value = 1.602267863950502
$ws.Cells.Item(8, 8).Value = 1.453362255853364
$ws.Cells.Item(8, 9).Value = 1.456683924343118
$ws.Cells.Item(8, 10).Value = 0.195091417285866
$ws.Cells.Item(8, 11).Value = 0.5157616409712773
$ws.Cells.Item(8, 12).Value = 0.3562565347776427
$ws.Cells.Item(8, 13).Value = 0.2576034550639612
$ws.Cells.Item(8, 14).Value = 2.706469710242349

$ws.Cells.Item(9, 2).Value = 0.8735754367883715
$ws.Cells.Item(9, 3).Value = 0.05707039175304374
$ws.Cells.Item(9, 4).Value = 0.09849164704082369
$ws.Cells.Item(9, 6).Value = 2.27920920784095
$ws.Cells.Item(9, 7).Value = 1.596082922521006
$ws.Cells.Item(9, 8).Value = 1.439837053660426
$ws.Cells.Item(9, 9).Value = 1.441719358917901
$ws.Cells.Item(9, 10).Value = 0.1929672116251364
$ws.Cells.Item(9, 11).Value = 0.5986438163095613
$ws.Cells.Item(9, 12).Value = 0.3613934634362721
$ws.Cells.Item(9, 13).Value = 0.2724195084811711
$ws.Cells.Item(9, 14).Value = 2.655283565420252

$ws.Cells.Item(10, 2).Value = 0.9349088831546339
$ws.Cells.Item(10, 3).Value = 0.06333294722767846
$ws.Cells.Item(10, 4).Value = 0.1003384891251713
$ws.Cells.Item(10, 6).Value = 2.276408862048029
$ws.Cells.Item(10, 7).Value = 1.595064208015202
$ws.Cells.Item(10, 8).Value = 1.432313752868382
$ws.Cells.Item(10, 9).Value = 1.433354223068797
$ws.Cells.Item(10, 10).Value = 0.1915880340351901
$ws.Cells.Item(10, 11).Value = 0.6605673120712368
$ws.Cells.Item(10, 12).Value = 0.3657915558167844
$ws.Cells.Item(10, 13).Value = 0.2838334619664735
$ws.Cells.Item(10, 14).Value = 2.621118274791762

$ws.Cells.Item(11, 2).Value = 0.9631659792015057
$ws.Cells.Item(11, 3).Value = 0.06615963675186265
$ws.Cells.Item(11, 4).Value = 0.1012076728927056
$ws.Cells.Item(11, 6).Value = 2.276099007185479
$ws.Cells.Item(11, 7).Value = 1.595365817147979
$ws.Cells.Item(11, 8).Value = 1.429413653935498
$ws.Cells.Item(11, 9).Value = 1.430118373289247
$ws.Cells.Item(11, 10).Value = 0.190999768458072
$ws.Cells.Item(11, 11).Value = 0.6889591861547331
$ws.Cells.Item(11, 12).Value = 0.3679268810836902
$ws.Cells.Item(11, 13).Value = 0.2891397011240002
$ws.Cells.Item(11, 14).Value = 2.60632146018029

$ws.Cells.Item(12, 2).Value = 0.9739169347422205
$ws.Cells.Item(12, 3).Value = 0.06722685799780947
$ws.Cells.Item(12, 4).Value = 0.1015409468064448
$ws.Cells.Item(12, 6).Value = 2.276120072990096
$ws.Cells.Item(12, 7).Value = 1.5955899708771
$ws.Cells.Item(12, 8).Value = 1.42839042539245
$ws.Cells.Item(12, 9).Value = 1.428974824534727
$ws.Cells.Item(12, 10).Value = 0.1907826155411154
$ws.Cells.Item(12, 11).Value = 0.6997421034867841
$ws.Cells.Item(12, 12).Value = 0.3687547336574539
$ws.Cells.Item(12, 13).Value = 0.2911653073978826
$ws.Cells.Item(12, 14).Value = 2.600825382267164

$ws.Cells.Item(13, 2).Value = 0.9715992831638118
$ws.Cells.Item(13, 3).Value = 0.06699715437973452
$ws.Cells.Item(13, 4).Value = 0.1014689871923196
$ws.Cells.Item(13, 6).Value = 2.276109384257538
$ws.Cells.Item(13, 7).Value = 1.595536806806621
$ws.Cells.Item(13, 8).Value = 1.428607463424868
$ws.Cells.Item(13, 9).Value = 1.429217471765682
$ws.Cells.Item(13, 10).Value = 0.1908291340364681
$ws.Cells.Item(13, 11).Value = 0.6974184145174434
$ws.Cells.Item(13, 12).Value = 0.3685755862702393
$ws.Cells.Item(13, 13).Value = 0.2907283364206137
$ws.Cells.Item(13, 14).Value = 2.60200429516448

$ws.Cells.Item(14, 2).Value = 0.9640494560777597
$ws.Cells.Item(14, 3).Value = 0.06624750148654357
$ws.Cells.Item(14, 4).Value = 0.1012350089592218
$ws.Cells.Item(14, 6).Value = 2.27609796780996
$ws.Cells.Item(14, 7).Value = 1.595382055425659
$ws.Cells.Item(14, 8).Value = 1.429327970409744
$ws.Cells.Item(14, 9).Value = 1.43002265404612
$ws.Cells.Item(14, 10).Value = 0.1909817908068394
$ws.Cells.Item(14, 11).Value = 0.6898456740452161
$ws.Cells.Item(14, 12).Value = 0.3679946037942159
$ws.Cells.Item(14, 13).Value = 0.2893060241296581
$ws.Cells.Item(14, 14).Value = 2.60586714697369

$ws.Cells.Item(15, 2).Value = 0.9594315411979721
$ws.Cells.Item(15, 3).Value = 0.06578790282904379
$ws.Cells.Item(15, 4).Value = 0.1010922274583024
$ws.Cells.Item(15, 6).Value = 2.276108992118765
$ws.Cells.Item(15, 7).Value = 1.595301581240989
$ws.Cells.Item(15, 8).Value = 1.429779062032537
$ws.Cells.Item(15, 9).Value = 1.430526501062815
$ws.Cells.Item(15, 10).Value = 0.1910760276584806
$ws.Cells.Item(15, 11).Value = 0.685211244016358
$ws.Cells.Item(15, 12).Value = 0.3676412391983632
$ws.Cells.Item(15, 13).Value = 0.2884369286631809
$ws.Cells.Item(15, 14).Value = 2.608247210095183

$ws.Cells.Item(16, 2).Value = 0.933069330620242
$ws.Cells.Item(16, 3).Value = 0.06314777052246257
$ws.Cells.Item(16, 4).Value = 0.1002822666476035
$ws.Cells.Item(16, 6).Value = 2.276448491319002
$ws.Cells.Item(16, 7).Value = 1.595059885540209
$ws.Cells.Item(16, 8).Value = 1.4325137799282
$ws.Cells.Item(16, 9).Value = 1.433577144967941
$ws.Cells.Item(16, 10).Value = 0.1916272645624879
$ws.Cells.Item(16, 11).Value = 0.6587162744203852
$ws.Cells.Item(16, 12).Value = 0.3656547069704175
$ws.Cells.Item(16, 13).Value = 0.2834889681854662
$ws.Cells.Item(16, 14).Value = 2.622100275911542

$ws.Cells.Item(17, 2).Value = 0.9169877629440748
$ws.Cells.Item(17, 3).Value = 0.06152245770553577
$ws.Cells.Item(17, 4).Value = 0.09979279003413666
$ws.Cells.Item(17, 6).Value = 2.276903514775569
$ws.Cells.Item(17, 7).Value = 1.595107505646908
$ws.Cells.Item(17, 8).Value = 1.434325119649841
$ws.Cells.Item(17, 9).Value = 1.43559440904005
$ws.Cells.Item(17, 10).Value = 0.1919754408829739
$ws.Cells.Item(17, 11).Value = 0.6425191360994233
$ws.Cells.Item(17, 12).Value = 0.3644704356665898
$ws.Cells.Item(17, 13).Value = 0.2804826461482151
$ws.Cells.Item(17, 14).Value = 2.630789534619982

$ws.Cells.Item(18, 2).Value = 0.9077716349054015
$ws.Cells.Item(18, 3).Value = 0.06058553297124547
$ws.Cells.Item(18, 4).Value = 0.09951399260103955
$ws.Cells.Item(18, 6).Value = 2.277255985368129
$ws.Cells.Item(18, 7).Value = 1.595206907580774
$ws.Cells.Item(18, 8).Value = 1.435416122207116
$ws.Cells.Item(18, 9).Value = 1.436808294929335
$ws.Cells.Item(18, 10).Value = 0.1921793866846713
$ws.Cells.Item(18, 11).Value = 0.6332239477910093
$ws.Cells.Item(18, 12).Value = 0.3638019506655326
$ws.Cells.Item(18, 13).Value = 0.2787642254525693
$ws.Cells.Item(18, 14).Value = 2.63585750279729

$ws.Cells.Item(19, 2).Value = 0.9046569957588133
$ws.Cells.Item(19, 3).Value = 0.06026794763445764
$ws.Cells.Item(19, 4).Value = 0.09942006792869762
$ws.Cells.Item(19, 6).Value = 2.277390919701517
$ws.Cells.Item(19, 7).Value = 1.595252933632295
$ws.Cells.Item(19, 8).Value = 1.435793965487306
$ws.Cells.Item(19, 9).Value = 1.437228506132996
$ws.Cells.Item(19, 10).Value = 0.1922490724555512
$ws.Cells.Item(19, 11).Value = 0.6300803740354013
$ws.Cells.Item(19, 12).Value = 0.3635777935585196
$ws.Cells.Item(19, 13).Value = 0.2781842459774495
$ws.Cells.Item(19, 14).Value = 2.637585475690436

$ws.Cells.Item(20, 2).Value = 0.9186962039478317
$ws.Cells.Item(20, 3).Value = 0.06169569107582618
$ws.Cells.Item(20, 4).Value = 0.09984461275010403
$ws.Cells.Item(20, 6).Value = 2.276845685781453
$ws.Cells.Item(20, 7).Value = 1.595094983922323
$ws.Cells.Item(20, 8).Value = 1.434127211686601
$ws.Cells.Item(20, 9).Value = 1.435374120199462
$ws.Cells.Item(20, 10).Value = 0.1919379957584564
$ws.Cells.Item(20, 11).Value = 0.6442411814879563
$ws.Cells.Item(20, 12).Value = 0.3645951921613033
$ws.Cells.Item(20, 13).Value = 0.2808015641218802
$ws.Cells.Item(20, 14).Value = 2.629857289324775

$ws.Cells.Item(21, 2).Value = 0.9662656547591837
$ws.Cells.Item(21, 3).Value = 0.06646777911066692
$ws.Cells.Item(21, 4).Value = 0.1013036222228365
$ws.Cells.Item(21, 6).Value = 2.276097566729092
$ws.Cells.Item(21, 7).Value = 1.595424526442187
$ws.Cells.Item(21, 8).Value = 1.42911430623468
$ws.Cells.Item(21, 9).Value = 1.429783933289542
$ws.Cells.Item(21, 10).Value = 0.1909367996747235
$ws.Cells.Item(21, 11).Value = 0.6920691195297479
$ws.Cells.Item(21, 12).Value = 0.3681647307992364
$ws.Cells.Item(21, 13).Value = 0.2897233518640689
$ws.Cells.Item(21, 14).Value = 2.604729625571469

$ws.Cells.Item(22, 2).Value = 0.9976497288341477
$ws.Cells.Item(22, 3).Value = 0.06956805343568817
$ws.Cells.Item(22, 4).Value = 0.1022812385931999
$ws.Cells.Item(22, 6).Value = 2.276415230943911
$ws.Cells.Item(22, 7).Value = 1.596280689394348
$ws.Cells.Item(22, 8).Value = 1.426275041235868
$ws.Cells.Item(22, 9).Value = 1.426607156897148
$ws.Cells.Item(22, 10).Value = 0.1903151549368296
$ws.Cells.Item(22, 11).Value = 0.7235110072341229
$ws.Cells.Item(22, 12).Value = 0.3706097926711891
$ws.Cells.Item(22, 13).Value = 0.2956489002896916
$ws.Cells.Item(22, 14).Value = 2.588931827796184

$ws.Cells.Item(23, 2).Value = 0.9808726971299961
$ws.Cells.Item(23, 3).Value = 0.06791507558098431
$ws.Cells.Item(23, 4).Value = 0.1017572784628413
$ws.Cells.Item(23, 6).Value = 2.276171957832929
$ws.Cells.Item(23, 7).Value = 1.595765129726885
$ws.Cells.Item(23, 8).Value = 1.427750470854591
$ws.Cells.Item(23, 9).Value = 1.428259072246021
$ws.Cells.Item(23, 10).Value = 0.1906439521927581
$ws.Cells.Item(23, 11).Value = 0.7067132404868062
$ws.Cells.Item(23, 12).Value = 0.3692945891703516
$ws.Cells.Item(23, 13).Value = 0.2924777100926335
$ws.Cells.Item(23, 14).Value = 2.59730625847693

$ws.Cells.Item(24, 2).Value = 0.9179237262096649
$ws.Cells.Item(24, 3).Value = 0.06161738008965756
$ws.Cells.Item(24, 4).Value = 0.09982117557015613
$ws.Cells.Item(24, 6).Value = 2.276871547191377
$ws.Cells.Item(24, 7).Value = 1.595100420634921
$ws.Cells.Item(24, 8).Value = 1.43421653124129
$ws.Cells.Item(24, 9).Value = 1.435473544151201
$ws.Cells.Item(24, 10).Value = 0.1919549129401652
$ws.Cells.Item(24, 11).Value = 0.6434625925202795
$ws.Cells.Item(24, 12).Value = 0.3645387512169549
$ws.Cells.Item(24, 13).Value = 0.2806573502902836
$ws.Cells.Item(24, 14).Value = 2.630278531902334

$ws.Cells.Item(25, 2).Value = 0.8513085153765303
$ws.Cells.Item(25, 3).Value = 0.05474443172397514
$ws.Cells.Item(25, 4).Value = 0.09783756880555217
$ws.Cells.Item(25, 6).Value = 2.281096702722849
$ws.Cells.Item(25, 7).Value = 1.597136819949853
$ws.Cells.Item(25, 8).Value = 1.443071529938464
$ws.Cells.Item(25, 9).Value = 1.445305480050109
$ws.Cells.Item(25, 10).Value = 0.1935099122984614
$ws.Cells.Item(25, 11).Value = 0.5760401749871562
$ws.Cells.Item(25, 12).Value = 0.359893860218051
$ws.Cells.Item(25, 13).Value = 0.2683181610409378
$ws.Cells.Item(25, 14).Value = 2.668526019250489
